# Insert a new daily price record for "Feria Lagunitas de Puerto Montt - Cebollín"
# as row 270, shifting the existing rows 270-344 down to 271-345.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 270 (pushes old 270..344 -> 271..345)
$ws.Rows.Item(270).Insert()

# Fill in the new row 270 with the new observation
$ws.Range("A270").Value = 4
$ws.Range("B270").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C270").Value = "Los Lagos"
$ws.Range("D270").Value = 44841
$ws.Range("E270").Value = 10
$ws.Range("F270").Value = 100112037
$ws.Range("G270").Value = "Cebollín"
$ws.Range("H270").Value = "Sin especificar"
$ws.Range("I270").Value = "Primera"
$ws.Range("J270").Value = 180
$ws.Range("K270").Value = 7000
$ws.Range("L270").Value = 7000
$ws.Range("M270").Value = 7000
$ws.Range("N270").Value = "$/paquete 36 unidades"
$ws.Range("O270").Value = "Región Metropolitana"
$ws.Range("P270").Value = 194
$ws.Range("Q270").Value = 36
$ws.Range("R270").Value = "Hortaliza"
